$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 06:28"

# Row 5: India - refreshed case counts
$ws.Range("B5").Value = 6623815
$ws.Range("C5").Value = 1635
$ws.Range("D5").Value = 5586703
$ws.Range("E5").Value = 934398

# Rows 52-55: Honduras jumps ahead of Costa Rica / Portugal / Etiopia
$ws.Range("A52").Value = "Honduras"
$ws.Range("B52").Value = 79629
$ws.Range("C52").Value = 841
$ws.Range("D52").Value = 29305
$ws.Range("E52").Value = 47902
$ws.Range("G52").Value = 23
$ws.Range("H52").Value = 2422

$ws.Range("A53").Value = "Costa Rica"
$ws.Range("B53").Value = 79182
$ws.Range("D53").Value = 45007
$ws.Range("E53").Value = 33225
$ws.Range("H53").Value = 950

$ws.Range("A54").Value = "Portugal"
$ws.Range("B54").Value = 79151
$ws.Range("D54").Value = 50207
$ws.Range("E54").Value = 26939
$ws.Range("H54").Value = 2005

$ws.Range("A55").Value = "Etiopia"
$ws.Range("B55").Value = 78819
$ws.Range("D55").Value = 33060
$ws.Range("E55").Value = 44537
$ws.Range("H55").Value = 1222

# Row 142: Tailandia - refreshed case counts
$ws.Range("B142").Value = 3590
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 3390
$ws.Range("E142").Value = 141

# Rows 153-155: Belice jumps ahead of Burkina Faso / Uruguay
$ws.Range("A153").Value = "Belice"
$ws.Range("B153").Value = 2196
$ws.Range("C153").Value = 65
$ws.Range("D153").Value = 1378
$ws.Range("E153").Value = 788
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 30

$ws.Range("A154").Value = "Burkina Faso"
$ws.Range("B154").Value = 2167
$ws.Range("D154").Value = 1419
$ws.Range("E154").Value = 689
$ws.Range("H154").Value = 59

$ws.Range("A155").Value = "Uruguay"
$ws.Range("B155").Value = 2145
$ws.Range("D155").Value = 1844
$ws.Range("E155").Value = 253
$ws.Range("H155").Value = 48

# Row 186: Mongolia - refreshed case counts
$ws.Range("B186").Value = 314
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 7

# Row 187: Butan - refreshed case counts
$ws.Range("B187").Value = 298
$ws.Range("C187").Value = 15
$ws.Range("D187").Value = 237
$ws.Range("E187").Value = 61

# Row 188: Camboya - refreshed case counts
$ws.Range("B188").Value = 280
$ws.Range("C188").Value = 2
$ws.Range("E188").Value = 5

# Rows 215-216: Montserrat swaps ahead of Islas Malvinas
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
